# Tried to implement Penality Reward System (unfinished)
# Updates the weekly forecast table (Sheet "Forecast Comparison") with a
# new week appended at the end and refreshed MyForecast numbers, then
# refreshes the dependent figures on the "Summary" sheet.

$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Sheet "Forecast Comparison": columns B (Week_Start_Date) and D (MyForecast) ---
# Force column B (rows 2-17) to stay text so the date-like strings are not
# auto-converted into Excel date serial numbers.
$wsForecast.Range("B2:B17").NumberFormat = "@"

$weekDates = @(
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20",
    "2025-04-27"
)

$myForecast = @(210, 199, 206, 229, 252, 165, 156, 154, 149, 216, 196, 159, 123, 106, 109, 109)

for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 2
    $wsForecast.Cells.Item($row, 2).Value = $weekDates[$i]
    $wsForecast.Cells.Item($row, 4).Value = $myForecast[$i]
}

# --- Sheet "Summary": refresh the derived metrics ---
$wsSummary.Cells.Item(2, 2).Value = "2023-01-08 to 2025-01-05"

# These cells hold plain numbers/dates but are stored as text on the sheet,
# so force a text number format before writing to avoid automatic
# conversion to numeric/date values.
$textCells = @(4, 5, 6, 7, 9, 10, 11, 12, 13, 14)
foreach ($r in $textCells) {
    $wsSummary.Cells.Item($r, 2).NumberFormat = "@"
}

$wsSummary.Cells.Item(4, 2).Value  = "376"
$wsSummary.Cells.Item(5, 2).Value  = "118"
$wsSummary.Cells.Item(6, 2).Value  = "88"
$wsSummary.Cells.Item(7, 2).Value  = "97"
$wsSummary.Cells.Item(8, 2).Value  = "12184 units"
$wsSummary.Cells.Item(9, 2).Value  = "2738"
$wsSummary.Cells.Item(10, 2).Value = "1571"
$wsSummary.Cells.Item(11, 2).Value = "844"
$wsSummary.Cells.Item(12, 2).Value = "252"
$wsSummary.Cells.Item(13, 2).Value = "2025-02-09"
$wsSummary.Cells.Item(14, 2).Value = "106"
